$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row: reorder L1/M1/N1 (Rappel columns moved before Caution) ---
$ws.Range("L1").Value = "MT brut (Rappel)"
$ws.Range("M1").Value = "Taxe (Rappel)"
$ws.Range("N1").Value = "Caution"

# --- Replace old row 2 (blank placeholder) with 6 data rows + a new total row 8 ---

# Row 2
$ws.Range("A2").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "B12346"
$ws.Range("D2").Value = "BAKKALI MOHAMED"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 10000

# Row 3
$ws.Range("A3").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "A123456"
$ws.Range("D3").Value = "YOUSSEF"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = "--"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 2000
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 18000

# Row 4
$ws.Range("A4").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "J207703"
$ws.Range("D4").Value = "ACHENGLI LAILA"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "--"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "--"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 10000

# Row 5
$ws.Range("A5").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "B12346"
$ws.Range("D5").Value = "BAKKALI MOHAMED"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 22000

# Row 6
$ws.Range("A6").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("B6").Value = "Direction régionale"
$ws.Range("C6").Value = "A123456"
$ws.Range("D6").Value = "YOUSSEF"
$ws.Range("E6").Value = "non"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 4000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 40000
$ws.Range("O6").Value = 43600

# Row 7
$ws.Range("A7").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("B7").Value = "Direction régionale"
$ws.Range("C7").Value = "J207703"
$ws.Range("D7").Value = "ACHENGLI LAILA"
$ws.Range("E7").Value = "non"
$ws.Range("F7").Value = "mensuelle"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 22000

# Row 8 (the former "total" row, now moved down with updated totals)
$ws.Range("A8").Value = " "
$ws.Range("B8").Value = " "
$ws.Range("C8").Value = " "
$ws.Range("D8").Value = " "
$ws.Range("E8").Value = " "
$ws.Range("F8").Value = " "
$ws.Range("G8").Value = " "
$ws.Range("H8").Value = 8000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 40000
$ws.Range("M8").Value = 2000
$ws.Range("N8").Value = 80000
$ws.Range("O8").Value = 125600

Write-Output "edit applied"
